$d = $word.ActiveDocument

# Locate the paragraph that starts the trailing "site footer" block, without
# relying on a hard-coded paragraph index.
$hit = $d.Content
$found = $hit.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx",
                            $false, $true, $false, $false, $false,
                            $true, 1, $false, "", 0)

if ($found) {
    $startIndex = -1
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $hit.Start -and $p.Range.End -ge $hit.End) {
            $startIndex = $i
            break
        }
    }

    if ($startIndex -gt 0) {
        # Remove this paragraph ("Ver no Jupiter Salvar em pdf Salvar em
        # docx"), the copyright/footer paragraph right after it, and the
        # blank paragraph that follows that -- three whole paragraphs
        # (including their paragraph marks). The blank paragraph that
        # already sits between the bibliography text and this footer block
        # is left in place, so exactly one blank line remains before the
        # page-break paragraph.
        $delStart = $d.Paragraphs.Item($startIndex).Range.Start
        $delEnd = $d.Paragraphs.Item($startIndex + 2).Range.End
        $d.Range($delStart, $delEnd).Delete()
    }
}
